# Apply the tracked changes to the "covariance" worksheet of the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("covariance")
[void]$ws.Activate()

# ---------------------------------------------------------------------------
# 1. Reword three of the "Notes" column explanations (column K).
# ---------------------------------------------------------------------------
$ws.Range("K6").Value = "Sum of deviation products"
$ws.Range("K7").Value = "Number of observations in the dataset"
$ws.Range("K8").Value = "Averaging deviation products by (N-1)"

# ---------------------------------------------------------------------------
# 2. Resize several columns.
# ---------------------------------------------------------------------------
# ColumnWidth values are offset by 5/6 from the stored OOXML "width" value,
# so subtract that offset to land on the target stored widths.
$ws.Columns.Item(1).ColumnWidth = 11.1640625 - (5/6)
$ws.Columns.Item(2).ColumnWidth = 11.1640625 - (5/6)
$ws.Columns.Item(3).ColumnWidth = 11.1640625 - (5/6)
$ws.Columns.Item(6).ColumnWidth = 11.83203125 - (5/6)
$ws.Columns.Item(7).ColumnWidth = 11.83203125 - (5/6)
$ws.Columns.Item(8).ColumnWidth = 24 - (5/6)

# ---------------------------------------------------------------------------
# 3. Turn on AutoFilter for the data table.
# ---------------------------------------------------------------------------
[void]$ws.Range("A5:H50").AutoFilter()

# ---------------------------------------------------------------------------
# 4. Update the active cell selection shown when the sheet is opened.
# ---------------------------------------------------------------------------
[void]$ws.Range("K16").Select()
